# ----------------------------------------------------------------------------
# edit.ps1
#
# Reproduces, through the PowerPoint COM object model, the two changes
# recorded in the source commit:
#
#  1. The table on slide 16 switches its table style (tableStyleId) from
#     {9695DEB9-E583-4B0E-BDD2-082FAC0E940B} to
#     {A05BE2D4-5C20-44A7-858A-677D192C8417}.
#
#  2. The presentation's theme colour scheme is swapped: the deck's main
#     design (used by every slide, via the slide master) moves away from
#     the custom "Integral" palette and over to the default "Office Theme"
#     palette (the colours that used to live only in the secondary theme
#     part used by the notes master).
# ----------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Table style change on slide 16.
# ---------------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
for ($i = 1; $i -le $slide16.Shapes.Count; $i++) {
    $shp = $slide16.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{A05BE2D4-5C20-44A7-858A-677D192C8417}")
    }
}

# ---------------------------------------------------------------------------
# 2) Theme colour swap: push the "Office Theme" palette onto the deck's
#    live design (slide master), replacing the "Integral" colours.
# ---------------------------------------------------------------------------
function HexToComRgb([string]$hex) {
    $rr = $hex.Substring(0, 2)
    $gg = $hex.Substring(2, 2)
    $bb = $hex.Substring(4, 2)
    # COM RGB longs are stored 0x00BBGGRR
    return [Convert]::ToInt32($bb + $gg + $rr, 16)
}

# Index -> (slot, target "Office Theme" colour) using the standard
# ThemeColorScheme ordering: 1=dk1 2=lt1 3=dk2 4=lt2 5-10=accent1-6
# 11=hlink 12=folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Item($i).RGB = (HexToComRgb $officeThemeColors[$i - 1])
}
